# Add the "git repository" call-out to the last slide (slide 16):
# a URL text box and, above it, a short label text box.
$p = $ppt.ActivePresentation
$s = $p.Slides.Item($p.Slides.Count)

# EMU -> point helper: PowerPoint's COM interface expresses shape
# geometry (Left/Top/Width/Height) in points, while the OOXML stores
# English Metric Units (1 pt = 12700 EMU).
$EMU = 12700

# --- Text box 1: hyperlink to the GitHub repository -----------------
$urlBox = $s.Shapes.AddTextbox(1, 6046879 / $EMU, 5870549 / $EMU, 5728224 / $EMU, 646331 / $EMU)
$urlBox.Name = "CasellaDiTesto 2"
$urlBox.Fill.Visible = 0
$urlBox.TextFrame.WordWrap = 1
$urlBox.TextFrame.AutoSize = 1

$urlRange = $urlBox.TextFrame.TextRange
$urlRange.Text = "https://github.com/alexsilza898/InsideSherpa/tree/master/KPMG/DataQ"
$urlRange.LanguageID = "it-IT"

# --- Text box 2: short caption above the link ------------------------
$labelBox = $s.Shapes.AddTextbox(1, 6046879 / $EMU, 5366747 / $EMU, 5186638 / $EMU, 369332 / $EMU)
$labelBox.Name = "CasellaDiTesto 5"
$labelBox.Fill.Visible = 0
$labelBox.TextFrame.WordWrap = 1
$labelBox.TextFrame.AutoSize = 1

$labelRange = $labelBox.TextFrame.TextRange
$labelRange.Text = "All code can be found at the git repository:"
$labelRange.LanguageID = "en-US"
